$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial number that was bumped
# from 45182 (2023-09-13) to 45184 (2023-09-15) for every data row
# (rows 2 through 141).
$ws.Range("C2:C141").Value = 45184
